$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new day column BL (16-aug) mirroring the style
# of the existing last column BK, then fill in the header + 24 hourly values.
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header cell's formatting (bold, border, centered/top alignment)
# from BK1 onto BL1 so the new column matches the rest of the header row.
$wsSpot.Range("BK1").Copy()
$wsSpot.Range("BL1").PasteSpecial(-4122)
$wsSpot.Range("BL1").Value = "16-aug"

$spotValues = @(
    76.70999999999999,
    74.77,
    69.38,
    51.61,
    44.57,
    45.13,
    38.97,
    48.95,
    54.8,
    60.93,
    26.28,
    3.8,
    0.17,
    -0.01,
    -0.01,
    -0.01,
    3.6,
    21.52,
    68.8,
    86.40000000000001,
    83.18000000000001,
    93.2,
    95.09999999999999,
    86.43000000000001
)

for ($i = 0; $i -lt $spotValues.Count; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 64).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 61 with the new daily price.
# The date column stores plain text dates ("2025-08-13", ...), not real
# Excel dates, so force text formatting before the write (otherwise "2025-
# 08-14" auto-converts to a date serial) and drop back to the default/
# unstyled cell format afterwards, matching the rest of the column.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A61").NumberFormat = "@"
$wsGaz.Range("A61").Value = "2025-08-14"
$wsGaz.Range("A61").Style = "Normal"
$wsGaz.Range("B61").Value = 31.325

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 61 with the new daily price.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A61").NumberFormat = "@"
$wsCo2.Range("A61").Value = "2025-08-14"
$wsCo2.Range("A61").Style = "Normal"
$wsCo2.Range("B61").Value = 70.48
